# Auto-generated Excel COM-interop script applying the Ragnarok_Profits scheduled-runner update.
# Updates cached market-price / profit columns (H:N) on each job sheet to match the refreshed data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1289
$ws.Range("I4").Value = 1016.125
$ws.Range("K4").Value = 1016.125
$ws.Range("M4").Value = -902.125
$ws.Range("H9").Value = 227.25
$ws.Range("I9").Value = 227.25
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 227.25
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -58.25
$ws.Range("N9").Value = ""
$ws.Range("H12").Value = 833.5
$ws.Range("I12").Value = 633.6667
$ws.Range("K12").Value = 633.6667
$ws.Range("M12").Value = -463.6667
$ws.Range("H33").Value = 239.94737
$ws.Range("I33").Value = 225.625
$ws.Range("J33").Value = 316.33334
$ws.Range("K33").Value = 225.625
$ws.Range("L33").Value = 316.33334
$ws.Range("M33").Value = 3.375
$ws.Range("N33").Value = -774.33334
$ws.Range("H43").Value = 3371.8333
$ws.Range("I43").Value = 4357.25
$ws.Range("J43").Value = 1401
$ws.Range("K43").Value = 4357.25
$ws.Range("L43").Value = 1401
$ws.Range("M43").Value = -4288.25
$ws.Range("N43").Value = -1539
$ws.Range("H45").Value = 1305
$ws.Range("I45").Value = 1425
$ws.Range("J45").Value = 1245
$ws.Range("K45").Value = 4275
$ws.Range("L45").Value = 3735
$ws.Range("M45").Value = -4083
$ws.Range("N45").Value = -4119
$ws.Range("H111").Value = 2543
$ws.Range("I111").Value = 2543
$ws.Range("K111").Value = 7629
$ws.Range("M111").Value = -4562
$ws.Range("H125").Value = 4010.5454
$ws.Range("J125").Value = 4021.6
$ws.Range("L125").Value = 36194.4
$ws.Range("N125").Value = -41114.4
$ws.Range("H135").Value = 3446.875
$ws.Range("I135").Value = 1471.2222
$ws.Range("J135").Value = 5987
$ws.Range("K135").Value = 13240.9998
$ws.Range("L135").Value = 53883
$ws.Range("M135").Value = -10705.9998
$ws.Range("N135").Value = -58953
$ws.Range("H138").Value = 6877.619
$ws.Range("I138").Value = 3654.7144
$ws.Range("K138").Value = 10964.1432
$ws.Range("M138").Value = -5824.143199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4449.25
$ws.Range("I63").Value = 4449.25
$ws.Range("K63").Value = 4449.25
$ws.Range("M63").Value = -3763.25
$ws.Range("H66").Value = 4449.25
$ws.Range("I66").Value = 4449.25
$ws.Range("K66").Value = 22246.25
$ws.Range("M66").Value = -18814.25
$ws.Range("H74").Value = 3076.2856
$ws.Range("I74").Value = 1907
$ws.Range("K74").Value = 1907
$ws.Range("M74").Value = -1033
$ws.Range("H77").Value = 3076.2856
$ws.Range("I77").Value = 1907
$ws.Range("K77").Value = 9535
$ws.Range("M77").Value = -5167
$ws.Range("H110").Value = 7809.25
$ws.Range("I110").Value = 8912.333000000001
$ws.Range("K110").Value = 8912.333000000001
$ws.Range("M110").Value = -6867.333000000001
$ws.Range("H122").Value = 3485.5557
$ws.Range("I122").Value = 4066.0667
$ws.Range("K122").Value = 12198.2001
$ws.Range("M122").Value = -9748.2001
$ws.Range("H132").Value = 6676209.5
$ws.Range("J132").Value = 20015006
$ws.Range("L132").Value = 60045018
$ws.Range("N132").Value = -60050078

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800
$ws.Range("H134").Value = 5885497.5
$ws.Range("I134").Value = 2896.3333
$ws.Range("K134").Value = 8688.999899999999
$ws.Range("M134").Value = -6153.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4764339.5
$ws.Range("I16").Value = 5265564.5
$ws.Range("K16").Value = 5265564.5
$ws.Range("M16").Value = -5265277.5
$ws.Range("H31").Value = 26318544
$ws.Range("I31").Value = 41669444
$ws.Range("K31").Value = 41669444
$ws.Range("M31").Value = -41669149
$ws.Range("H34").Value = 26318544
$ws.Range("I34").Value = 41669444
$ws.Range("K34").Value = 41669444
$ws.Range("M34").Value = -41669242
$ws.Range("H113").Value = 4764339.5
$ws.Range("I113").Value = 5265564.5
$ws.Range("K113").Value = 5265564.5
$ws.Range("M113").Value = -5263394.5
$ws.Range("H119").Value = 80633
$ws.Range("J119").Value = 80633
$ws.Range("L119").Value = 80633
$ws.Range("N119").Value = -90309
$ws.Range("H132").Value = 2871.8076
$ws.Range("I132").Value = 2670.7368
$ws.Range("J132").Value = 3417.5715
$ws.Range("K132").Value = 8012.2104
$ws.Range("L132").Value = 10252.7145
$ws.Range("M132").Value = -5482.2104
$ws.Range("N132").Value = -15312.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 135.83333
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2927.5557
$ws.Range("I80").Value = 1679.25
$ws.Range("J80").Value = 3926.2
$ws.Range("K80").Value = 1679.25
$ws.Range("L80").Value = 3926.2
$ws.Range("M80").Value = -681.25
$ws.Range("N80").Value = -5922.2
$ws.Range("H83").Value = 2927.5557
$ws.Range("I83").Value = 1679.25
$ws.Range("J83").Value = 3926.2
$ws.Range("K83").Value = 8396.25
$ws.Range("L83").Value = 19631
$ws.Range("M83").Value = -3404.25
$ws.Range("N83").Value = -29615
$ws.Range("H97").Value = 522.5909
$ws.Range("I97").Value = 631.64703
$ws.Range("K97").Value = 631.64703
$ws.Range("M97").Value = -135.64703
$ws.Range("H102").Value = 2764.543
$ws.Range("I102").Value = 2773.3447
$ws.Range("K102").Value = 2773.3447
$ws.Range("M102").Value = -1151.3447
$ws.Range("H132").Value = 50002500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 50002500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 150007500
$ws.Range("N132").Value = -150012560
$ws.Range("M132").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 485
$ws.Range("I22").Value = 485
$ws.Range("K22").Value = 485
$ws.Range("M22").Value = -190
$ws.Range("H27").Value = 485
$ws.Range("I27").Value = 485
$ws.Range("K27").Value = 485
$ws.Range("M27").Value = -378
$ws.Range("H68").Value = 4169163.2
$ws.Range("I68").Value = 10418516
$ws.Range("J68").Value = 2928.3333
$ws.Range("K68").Value = 10418516
$ws.Range("L68").Value = 2928.3333
$ws.Range("M68").Value = -10417767
$ws.Range("N68").Value = -4426.3333
$ws.Range("H71").Value = 4169163.2
$ws.Range("I71").Value = 10418516
$ws.Range("J71").Value = 2928.3333
$ws.Range("K71").Value = 52092580
$ws.Range("L71").Value = 14641.6665
$ws.Range("M71").Value = -52088836
$ws.Range("N71").Value = -22129.6665
$ws.Range("H82").Value = 6466.5
$ws.Range("I82").Value = 3380.7144
$ws.Range("J82").Value = 13666.667
$ws.Range("K82").Value = 3380.7144
$ws.Range("L82").Value = 13666.667
$ws.Range("M82").Value = -3019.7144
$ws.Range("N82").Value = -14388.667
$ws.Range("H85").Value = 6466.5
$ws.Range("I85").Value = 3380.7144
$ws.Range("J85").Value = 13666.667
$ws.Range("K85").Value = 3380.7144
$ws.Range("L85").Value = 13666.667
$ws.Range("M85").Value = -2132.7144
$ws.Range("N85").Value = -16162.667
$ws.Range("H93").Value = 2061301.5
$ws.Range("I93").Value = 1849.8
$ws.Range("K93").Value = 1849.8
$ws.Range("M93").Value = -601.8
$ws.Range("H100").Value = 22730822
$ws.Range("I100").Value = 3577.5715
$ws.Range("K100").Value = 3577.5715
$ws.Range("M100").Value = -3036.5715
$ws.Range("H103").Value = 48090.168
$ws.Range("J103").Value = 48090.168
$ws.Range("L103").Value = 48090.168
$ws.Range("N103").Value = -50434.168
$ws.Range("H132").Value = 4180.25
$ws.Range("I132").Value = 2742.6667
$ws.Range("J132").Value = 6336.625
$ws.Range("K132").Value = 8228.000100000001
$ws.Range("L132").Value = 19009.875
$ws.Range("M132").Value = -5698.000100000001
$ws.Range("N132").Value = -24069.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 99999
$ws.Range("I2").Value = 99999
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 99999
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -99887
$ws.Range("N2").Value = ""
$ws.Range("H81").Value = 1989
$ws.Range("I81").Value = 1794.7142
$ws.Range("K81").Value = 3589.4284
$ws.Range("M81").Value = -2528.4284
$ws.Range("H84").Value = 1989
$ws.Range("I84").Value = 1794.7142
$ws.Range("K84").Value = 17947.142
$ws.Range("M84").Value = -12643.142
$ws.Range("H122").Value = 2530.2942
$ws.Range("I122").Value = 2608.2144
$ws.Range("K122").Value = 7824.6432
$ws.Range("M122").Value = -5374.6432
$ws.Range("H132").Value = 631210
$ws.Range("I132").Value = 5769.5454
$ws.Range("K132").Value = 17308.6362
$ws.Range("M132").Value = -14778.6362
